$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so Excel does not auto-convert them (e.g. "1.00" -> 1, "0.0731" -> 7.31E-02).
# We do this by temporarily setting NumberFormat to "@" (text), assigning the
# value, then resetting the style back to "Normal" so no extra style lingers.
$textForceCells = @("D24", "D31", "D44", "D23", "D34", "D47", "D30", "D35", "D25", "D20", "D6", "D22", "D33", "D43", "D41", "D29", "D21", "D15", "D5", "D28", "D39", "D45", "D16", "D27", "D37")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply all updated values from the data refresh
$ws.Range("D2").Value = "64.111.73"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "3.417.55"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "571.71"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "160.78"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D8").Value = "3.419.09"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -8.56%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("D13").Value = "4.007.70"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "27.18"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "0.0000176"
$ws.Range("E16").Value = "  -6.66%  "
$ws.Range("D17").Value = "64.169.57"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "3.449.71"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  -3.81%  "
$ws.Range("D20").Value = "13.59"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").Value = "378.58"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "7.87"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "71.42"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "0.518"
$ws.Range("E25").Value = "  -5.45%  "
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  -4.94%  "
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "6.07"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -4.55%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "22.99"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").Value = "7.08"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  -4.75%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "0.859"
$ws.Range("E37").Value = "  +11.41%  "
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0731"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.808.78"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").Value = "43.12"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").Value = "25.74"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.45"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "26.15"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("D47").Value = "342.05"
$ws.Range("E47").Value = "  +7.69%  "
$ws.Range("E48").Value = "  +5.20%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("E51").Value = "  -4.45%  "

# Restore default style on the text-forced cells (removes the temporary
# text NumberFormat so styling matches the original unstyled cells)
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
